# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (copied from the existing "2022-Q2"
# worksheet so it inherits the same layout/formatting) right after the
# "总计" summary sheet, fills it with the new quarter's fund data, and
# prepends a matching summary row to "总计" (shifting the older rows
# down by one).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计" by duplicating
#    "2022-Q2" (keeps header/style/column layout identical) and then
#    overwrite its data with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# 2022-Q3 only has two funds on record, so drop the extra rows that
# were copied over from 2022-Q2 (rows 4 and 5).
$wsQ3.Rows.Item(5).Delete()
$wsQ3.Rows.Item(4).Delete()

# Row 2: 基金代码 012349 / 天弘恒生科技指数（QDII）C
$wsQ3.Range("B2:G2").NumberFormat = "@"
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "012349"
$wsQ3.Range("C2").Value = "天弘恒生科技指数（QDII）C"
$wsQ3.Range("D2").Value = "33.57"
$wsQ3.Range("E2").Value = "92.84"
$wsQ3.Range("F2").Value = "6.98"
$wsQ3.Range("G2").Value = "2.3432"
$wsQ3.Range("H2").Value = 5

# Row 3: 基金代码 012348 / 天弘恒生科技指数（QDII）A
$wsQ3.Range("B3:G3").NumberFormat = "@"
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "012348"
$wsQ3.Range("C3").Value = "天弘恒生科技指数（QDII）A"
$wsQ3.Range("D3").Value = "30.64"
$wsQ3.Range("E3").Value = "92.84"
$wsQ3.Range("F3").Value = "6.98"
$wsQ3.Range("G3").Value = "2.1387"
$wsQ3.Range("H3").Value = 5

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q3" row to the "总计" summary sheet. Insert a row
#    at position 2, copy the formatting from the row right below it
#    (which holds the existing data format) and fill in the new values.
#    The existing rows below keep their values; only the running index
#    in column A needs to be bumped by one.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 4.48

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
